# "Updated sine and cosine estimation."
#
# Bumps the per-unit quantities for the three top-of-BOM line items (rows
# 16-18) and widens the row-19 rollup formula's multipliers to match -
# every other changed cell in the workbook is a downstream formula result
# that recalculates automatically once these inputs change. Also re-points
# the frozen-pane scroll position and the current selection to match the
# saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Input edits -----------------------------------------------------
# Quantity column (C) for the three components feeding row 19's totals.
$ws.Range("C16").Value = 3
$ws.Range("C17").Value = 7
$ws.Range("C18").Value = 10

# Row 19's F column total used to scale 2x -> 3x for F4/F5 (e.g. extra
# fastener/hardware count per assembly).
$ws.Range("F19").Formula = "=3*F4+3*F5+F3"

# --- View state --------------------------------------------------------
# Re-anchor the frozen pane's scrolled-to row and the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$null = $ws.Range("B24:O24").Select()

# Best-effort: restore the saved workbook window size (not all hosts
# expose window chrome sizing through automation).
try {
    $win.Width = 25600
    $win.Height = 16060
} catch {
}
